$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update odds values in columns C, D, E for rows 1-9
$values = @{
    "C1" = 2.65;  "D1" = 3.35
    "C2" = 1.4;   "D2" = 4.75;  "E2" = 6.5
    "C3" = 2;     "D3" = 3.4;  "E3" = 3.5
    "D4" = 5.5;   "E4" = 11
    "C5" = 2.1;   "D5" = 3.4;  "E5" = 3.3
    "C6" = 2.15;  "D6" = 3.3;  "E6" = 3.2
    "C7" = 4.5;   "D7" = 3.75; "E7" = 1.69
    "C8" = 2.0499999999999998; "D8" = 3.4; "E8" = 3.4
    "C9" = 1.61;  "D9" = 3.75; "E9" = 5.25
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Update the selected cell on the active sheet (F9 -> F8)
$ws.Range("F8").Select()
